$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.582.34"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.925.03"
$ws.Range("E3").Value = "  +3.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.95%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.60"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4748"
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2924"
$ws.Range("E8").Value = "  +3.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06794"
$ws.Range("E9").Value = "  +6.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "106.03"
$ws.Range("E10").Value = "  +11.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.42"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").Value = "1.913.42"
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07725"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.355"
$ws.Range("E14").Value = "  +7.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6727"
$ws.Range("E15").Value = "  +5.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "288.94"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "30.610.91"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007650"
$ws.Range("E18").Value = "  +4.16%  "
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9995"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.512"
$ws.Range("E21").Value = "  +9.91%  "
$ws.Range("D22").Value = "2.159.82"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.300"
$ws.Range("E24").Value = "  +5.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.423"
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.39"
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.77"
$ws.Range("E27").Value = "  +7.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.134"
$ws.Range("E28").Value = "  +11.43%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.185"
$ws.Range("E31").Value = "  +4.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.070"
$ws.Range("E32").Value = "  +7.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05078"
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7427"
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.159"
$ws.Range("E35").Value = "  +3.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02085"
$ws.Range("E36").Value = "  +8.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.750"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.069"
$ws.Range("E39").Value = "  +4.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "111.03"
$ws.Range("E40").Value = "  +5.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8806"
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4414"
$ws.Range("E42").Value = "  +9.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.948"
$ws.Range("E43").Value = "  +6.93%  "
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("E45").Value = "  +3.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.286"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.404"
$ws.Range("E47").Value = "  +6.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1234"
$ws.Range("E48").Value = "  +4.09%  "
$ws.Range("B49").Value = "WOONetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.2545"
$ws.Range("E49").Value = "  +17.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "47.23"
$ws.Range("E50").Value = "  +16.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.30"
$ws.Range("E51").Value = "  +4.30%  "
